$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.73319999999998
$ws.Range("A6").Value = -22.66590000000001
$ws.Range("A7").Value = -21.99720000000001
$ws.Range("B7").Value = 5.076900000000003
$ws.Range("B12").Value = 4.804499999999999
$ws.Range("D13").Value = -8.377700000000001
$ws.Range("D14").Value = -8.323999999999998
$ws.Range("B15").Value = 5.193099999999998
$ws.Range("A16").Value = -21.55849999999999
$ws.Range("D16").Value = -8.816100000000004
$ws.Range("D19").Value = -7.261299999999997
$ws.Range("A20").Value = -23.1615
$ws.Range("B20").Value = 5.052599999999997
$ws.Range("B21").Value = 10.60020000000001
$ws.Range("B22").Value = 8.858700000000004
$ws.Range("D22").Value = -8.074800000000002
$ws.Range("B23").Value = 8.880500000000007
$ws.Range("A28").Value = -22.2095
$ws.Range("A29").Value = -21.6891
$ws.Range("B29").Value = 5.101500000000002
$ws.Range("A32").Value = -21.08599999999999
$ws.Range("B34").Value = 10.23770000000001
$ws.Range("D36").Value = -7.800999999999997
$ws.Range("A40").Value = -19.40009999999999
$ws.Range("B42").Value = 9.283300000000001
$ws.Range("B43").Value = 6.001400000000001
$ws.Range("B44").Value = 5.548999999999999
$ws.Range("B45").Value = 5.191900000000003
$ws.Range("A46").Value = -22.1632
$ws.Range("B46").Value = 5.530400000000006
$ws.Range("D46").Value = -7.940199999999999
$ws.Range("B50").Value = 4.911899999999993
$ws.Range("D50").Value = -7.923399999999998
$ws.Range("A51").Value = -22.308
$ws.Range("B51").Value = 5.620299999999999
$ws.Range("A52").Value = -22.0768
$ws.Range("A57").Value = -22.64180000000002
$ws.Range("A59").Value = -21.9102
$ws.Range("A62").Value = -21.9641
$ws.Range("A66").Value = -21.56649999999999
$ws.Range("B66").Value = 5.061899999999998
$ws.Range("B67").Value = 5.241000000000001
$ws.Range("A73").Value = -20.0855
$ws.Range("A74").Value = -21.92629999999998
$ws.Range("B79").Value = 9.877100000000002
$ws.Range("B84").Value = 5.501500000000002
$ws.Range("A92").Value = -21.41210000000001
$ws.Range("B92").Value = 5.059899999999995
$ws.Range("D95").Value = -8.307599999999999
$ws.Range("B97").Value = 5.827700000000001
$ws.Range("D97").Value = -8.667899999999996
$ws.Range("A100").Value = -21.989
